# The sheet originally has 4 columns: Nombre | NIF | Email (mailto: hyperlinks) | Colegio
# This edit drops the "Email" column (column C) entirely, which shifts the
# "Colegio" column (old D) left into column C, and removes the now-orphaned
# hyperlinks plus the "Hipervinculo" cell style that only that column used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "Email" column (C) - shifts "Colegio" (D) left to C,
# drops its shared-strings (emails) and its hyperlink-styled cell format.
$ws.Range("C1:C4").EntireColumn.Delete()

# Drop the now-orphaned hyperlink definitions that used to live on column C.
$ws.Hyperlinks.Delete()

# Remove the now-unused "Hipervinculo" named cell style (item 1 - the only
# non-"Normal" style in this workbook). Addressed by index rather than by
# its accented name to sidestep any encoding ambiguity in the script text.
$wb.Styles.Item(1).Delete()

# Match the post-edit selection: C1:C4 (the old "Email" hyperlink column's
# worth of rows, now holding "Colegio").
$ws.Range("C1:C4").Select()
